# TC06_Trials_Filter_Race-White.xlsx — add the missing Neo4j query text to
# the "startup" sheet (cell A2), matching the commit's "added queries in
# all ctdc tc xls" change, then leave the selection the way the author's
# Excel session left it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# NB: single-quoted PowerShell literal so the embedded backticks (used
# around `Case ID`, `Trial Code`, `Arm`, `Arm Treatment`) stay literal;
# embedded single quotes are doubled per PowerShell literal-string rules.
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN [''WHITE''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Range("A2").Value = $query
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 87

# Selection ends on B2:B5, mirroring the saved sheetView.
$ws.Range("B2:B5").Select()

$wb.Save()
